$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set B11 first so its string becomes sharedString index 9 (matching target
# insertion order), then B9 (index 10), then B10 (index 11).
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Implementation of a complex form to publish lessons timetable"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Implementation of a workflow allowing to save ""connections"" with other people"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Implementation of the Search filter by teacher's name"

$ws.Range("B10").Select()
